# Insert a new weekly price-report row for "Ají" (chili pepper) at row 233
# of the "Feria Lagunitas de Puerto Montt" sheet. Inserting the row pushes
# the existing rows 233-336 down to 234-337 (dimension grows from
# A1:R336 to A1:R337), and the new row 233 is populated with this week's
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 233:336 down to 234:337, leaving row 233 blank and ready
# for the new record.
$ws.Rows.Item(233).Insert()

# Populate the newly inserted row 233 with the new weekly record.
$ws.Cells.Item(233, 1).Value = 4
$ws.Cells.Item(233, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(233, 3).Value = "Los Lagos"
$ws.Cells.Item(233, 4).Value = 44875
$ws.Cells.Item(233, 5).Value = 10
$ws.Cells.Item(233, 6).Value = 100112021
$ws.Cells.Item(233, 7).Value = "Ají"
$ws.Cells.Item(233, 8).Value = "Inferno"
$ws.Cells.Item(233, 9).Value = "Primera"
$ws.Cells.Item(233, 10).Value = 80
$ws.Cells.Item(233, 11).Value = 21000
$ws.Cells.Item(233, 12).Value = 21000
$ws.Cells.Item(233, 13).Value = 21000
$ws.Cells.Item(233, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(233, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(233, 16).Value = 2100
$ws.Cells.Item(233, 17).Value = 10
$ws.Cells.Item(233, 18).Value = "Hortaliza"

# Match the date-number formatting used by the rest of column D.
$ws.Cells.Item(233, 4).NumberFormat = $ws.Cells.Item(234, 4).NumberFormat
